$wb = $excel.ActiveWorkbook

# --- Sheet "A3": ClassID values updated to the full IBC wasm contract path,
#     and the now-redundant helper row (class_id / wasm contract lookup) removed.
$wsA3 = $wb.Worksheets.Item("A3")
$wsA3.Range("B2").Value = "wasm.stars1ve46fjrhcrum94c7d8yc2wsdz8cpuw73503e8qn9r44spr6dw0lsvmvtqh/channel-207/ddnft"
$wsA3.Range("B3").Value = "wasm.stars1ve46fjrhcrum94c7d8yc2wsdz8cpuw73503e8qn9r44spr6dw0lsvmvtqh/channel-207/ddnft"
$wsA3.Rows.Item(5).Delete()
$wsA3.Range("E22").Select()

# --- Sheet "A4": ClassID values updated to the ibc/ denom hash, and the
#     now-redundant helper row (denom_id lookup) cleared out.
$wsA4 = $wb.Worksheets.Item("A4")
$wsA4.Range("B2").Value = "ibc/80F44B622197B3F4CF25D1D7510A7F5BA50B3BBFC9DCA9B7659FB6572103BFB5"
$wsA4.Range("B3").Value = "ibc/80F44B622197B3F4CF25D1D7510A7F5BA50B3BBFC9DCA9B7659FB6572103BFB5"
$wsA4.Range("E5:F5").ClearContents()
$wsA4.Range("E5:I5").Select()

# --- Sheet "A5": no data change, just cursor position.
$wsA5 = $wb.Worksheets.Item("A5")
$wsA5.Range("D2").Select()

# --- Sheet "A6": ChainID corrected from gon-irishub-1 to gon-flixnet-1.
#     This becomes the active tab/sheet at the end of the session.
$wsA6 = $wb.Worksheets.Item("A6")
$wsA6.Range("D2").Value = "gon-flixnet-1"
$wsA6.Range("D2").Select()
